$wb = $excel.ActiveWorkbook

# Sheet1 - weibull
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = -2.02920667547382
$ws.Range("C2").Value = 0.101386015116349
$ws.Range("B3").Value = -0.165252558924571
$ws.Range("C3").Value = 0.04858091302945

# Sheet2 - lognormal
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = 1.12639172344645
$ws.Range("C2").Value = 0.105907986747605
$ws.Range("B3").Value = -0.699866140753643
$ws.Range("C3").Value = 0.0462881852377329

# Sheet3 - llogis
$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = -1.61323348320391
$ws.Range("C2").Value = 0.121792974479335
$ws.Range("B3").Value = 0.180586724304846
$ws.Range("C3").Value = 0.0673717463977423

# Sheet4 - gompertz
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = -1.84249044736258
$ws.Range("C2").Value = 0.127874783169237
$ws.Range("B3").Value = -0.0639210532659732
$ws.Range("C3").Value = 0.0168349599084398

# Sheet5 - exp (no value changes)

# Sheet6 - weibull cov
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = 0.0102791240611726
$ws.Range("B2").Value = -0.0018486256997203
$ws.Range("A3").Value = -0.0018486256997203
$ws.Range("B3").Value = 0.00236010511077498

# Sheet7 - lognormal cov
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = 0.0112165016569309
$ws.Range("B2").Value = -0.00303414204430054
$ws.Range("A3").Value = -0.00303414204430054
$ws.Range("B3").Value = 0.00214259609260268

# Sheet8 - llogis cov
$ws = $wb.Worksheets.Item(8)
$ws.Range("A2").Value = 0.0148335286325239
$ws.Range("B2").Value = 0.00101856287708708
$ws.Range("A3").Value = 0.00101856287708708
$ws.Range("B3").Value = 0.0045389522126817

# Sheet9 - gompertz cov
$ws = $wb.Worksheets.Item(9)
$ws.Range("A2").Value = 0.0163519601705793
$ws.Range("B2").Value = -0.00153860905264522
$ws.Range("A3").Value = -0.00153860905264522
$ws.Range("B3").Value = 0.000283415875118776

# Sheet10 - exp cov (no value changes)
